$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 13889982
$ws.Range("J32").Value = 1585.3334
$ws.Range("L32").Value = 1585.3334
$ws.Range("N32").Value = -2237.3334
$ws.Range("H92").Value = 1093.7222
$ws.Range("I92").Value = 1178.3846
$ws.Range("J92").Value = 873.6
$ws.Range("K92").Value = 1178.3846
$ws.Range("L92").Value = 873.6
$ws.Range("M92").Value = 69.61539999999991
$ws.Range("N92").Value = -3369.6
$ws.Range("H112").Value = 1677.091
$ws.Range("J112").Value = 1950.6111
$ws.Range("L112").Value = 5851.8333
$ws.Range("N112").Value = -8067.8333
$ws.Range("H132").Value = 8075.684
$ws.Range("I132").Value = 8075.684
$ws.Range("K132").Value = 24227.052
$ws.Range("M132").Value = -21697.052
$ws.Range("H137").Value = 2787.6296
$ws.Range("I137").Value = 1713.8182
$ws.Range("K137").Value = 5141.4546
$ws.Range("M137").Value = -2591.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 880.8
$ws.Range("I2").Value = 829.9167
$ws.Range("J2").Value = 1084.3334
$ws.Range("K2").Value = 829.9167
$ws.Range("L2").Value = 1084.3334
$ws.Range("M2").Value = -716.9167
$ws.Range("N2").Value = -1310.3334
$ws.Range("H32").Value = 2066936
$ws.Range("I32").Value = 1060523.8
$ws.Range("J32").Value = 7937674
$ws.Range("K32").Value = 1060523.8
$ws.Range("L32").Value = 7937674
$ws.Range("M32").Value = -1060236.8
$ws.Range("N32").Value = -7938248
$ws.Range("H45").Value = 30974.934
$ws.Range("I45").Value = 63016.715
$ws.Range("K45").Value = 63016.715
$ws.Range("M45").Value = -62639.715
$ws.Range("H61").Value = 3065.6
$ws.Range("I61").Value = 2640.5715
$ws.Range("J61").Value = 3437.5
$ws.Range("K61").Value = 2640.5715
$ws.Range("L61").Value = 3437.5
$ws.Range("M61").Value = -2428.5715
$ws.Range("N61").Value = -3861.5
$ws.Range("H116").Value = 880.8
$ws.Range("I116").Value = 829.9167
$ws.Range("J116").Value = 1084.3334
$ws.Range("K116").Value = 829.9167
$ws.Range("L116").Value = 1084.3334
$ws.Range("M116").Value = 1464.0833
$ws.Range("N116").Value = -5672.3334
$ws.Range("H136").Value = 3065.6
$ws.Range("I136").Value = 2640.5715
$ws.Range("J136").Value = 3437.5
$ws.Range("K136").Value = 7921.7145
$ws.Range("L136").Value = 10312.5
$ws.Range("M136").Value = -5371.7145
$ws.Range("N136").Value = -15412.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 880.8
$ws.Range("I3").Value = 829.9167
$ws.Range("J3").Value = 1084.3334
$ws.Range("K3").Value = 829.9167
$ws.Range("L3").Value = 1084.3334
$ws.Range("M3").Value = -715.9167
$ws.Range("N3").Value = -1312.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1349
$ws.Range("I16").Value = 1199
$ws.Range("K16").Value = 1199
$ws.Range("M16").Value = -912
$ws.Range("H31").Value = 3792020
$ws.Range("I31").Value = 2811.3157
$ws.Range("J31").Value = 8934518
$ws.Range("K31").Value = 2811.3157
$ws.Range("L31").Value = 8934518
$ws.Range("M31").Value = -2516.3157
$ws.Range("N31").Value = -8935108
$ws.Range("H34").Value = 3792020
$ws.Range("I34").Value = 2811.3157
$ws.Range("J34").Value = 8934518
$ws.Range("K34").Value = 2811.3157
$ws.Range("L34").Value = 8934518
$ws.Range("M34").Value = -2609.3157
$ws.Range("N34").Value = -8934922
$ws.Range("H58").Value = 2199.3572
$ws.Range("I58").Value = 1804.65
$ws.Range("J58").Value = 3186.125
$ws.Range("K58").Value = 1804.65
$ws.Range("L58").Value = 3186.125
$ws.Range("M58").Value = -1601.65
$ws.Range("N58").Value = -3592.125
$ws.Range("H99").Value = 3563.75
$ws.Range("J99").Value = 3806.25
$ws.Range("L99").Value = 3806.25
$ws.Range("N99").Value = -6802.25
$ws.Range("H107").Value = 8334666.5
$ws.Range("J107").Value = 1865.6666
$ws.Range("L107").Value = 1865.6666
$ws.Range("N107").Value = -5705.6666
$ws.Range("H113").Value = 1349
$ws.Range("I113").Value = 1199
$ws.Range("K113").Value = 1199
$ws.Range("M113").Value = 971
$ws.Range("H122").Value = 1613.625
$ws.Range("I122").Value = 1604.25
$ws.Range("J122").Value = 1623
$ws.Range("K122").Value = 4812.75
$ws.Range("L122").Value = 4869
$ws.Range("M122").Value = -2362.75
$ws.Range("N122").Value = -9769
$ws.Range("H126").Value = 3563.75
$ws.Range("J126").Value = 3806.25
$ws.Range("L126").Value = 11418.75
$ws.Range("N126").Value = -16358.75
$ws.Range("H132").Value = 12351113
$ws.Range("I132").Value = 4327.533
$ws.Range("J132").Value = 27784594
$ws.Range("K132").Value = 12982.599
$ws.Range("L132").Value = 83353782
$ws.Range("M132").Value = -10452.599
$ws.Range("N132").Value = -83358842
$ws.Range("H134").Value = 3941.682
$ws.Range("I134").Value = 3941.682
$ws.Range("K134").Value = 11825.046
$ws.Range("M134").Value = -9290.045999999998
$ws.Range("H136").Value = 2199.3572
$ws.Range("I136").Value = 1804.65
$ws.Range("J136").Value = 3186.125
$ws.Range("K136").Value = 5413.950000000001
$ws.Range("L136").Value = 9558.375
$ws.Range("M136").Value = -2863.950000000001
$ws.Range("N136").Value = -14658.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2432
$ws.Range("I81").Value = 3198
$ws.Range("K81").Value = 9594
$ws.Range("M81").Value = -8471
$ws.Range("H84").Value = 2432
$ws.Range("I84").Value = 3198
$ws.Range("K84").Value = 28782
$ws.Range("M84").Value = -23166
$ws.Range("H107").Value = 440.05264
$ws.Range("I107").Value = 310.33334
$ws.Range("J107").Value = 480.31033
$ws.Range("K107").Value = 931.0000200000001
$ws.Range("L107").Value = 1440.93099
$ws.Range("M107").Value = 988.9999799999999
$ws.Range("N107").Value = -5280.93099
$ws.Range("H125").Value = 4999
$ws.Range("I125").Value = 4999
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 14997
$ws.Range("N125").ClearContents()
$ws.Range("M125").Value = -10077
$ws.Range("H131").Value = 6946557
$ws.Range("J131").Value = 1940.2424
$ws.Range("L131").Value = 5820.7272
$ws.Range("N131").Value = -15900.7272
$ws.Range("H132").Value = 1399.5
$ws.Range("J132").Value = 999
$ws.Range("L132").Value = 8991
$ws.Range("N132").Value = -14051
$ws.Range("H134").Value = 4052.4167
$ws.Range("I134").Value = 1255.8334
$ws.Range("J134").Value = 4984.6113
$ws.Range("K134").Value = 3767.5002
$ws.Range("L134").Value = 14953.8339
$ws.Range("M134").Value = 1302.4998
$ws.Range("N134").Value = -25093.8339
$ws.Range("H139").Value = 5171.4688
$ws.Range("I139").Value = 3006.6924
$ws.Range("K139").Value = 9020.0772
$ws.Range("M139").Value = -3880.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 525138.2
$ws.Range("I14").Value = 525138.2
$ws.Range("K14").Value = 525138.2
$ws.Range("M14").Value = -524970.2
$ws.Range("H132").Value = 62503030
$ws.Range("I132").Value = 200002400
$ws.Range("K132").Value = 600007200
$ws.Range("M132").Value = -600004670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 24503
$ws.Range("I11").Value = 24503
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 24503
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H93").Value = 3830
$ws.Range("I93").Value = 250
$ws.Range("J93").Value = 4725
$ws.Range("K93").Value = 250
$ws.Range("L93").Value = 4725
$ws.Range("M93").Value = 998
$ws.Range("N93").Value = -7221
$ws.Range("H100").Value = 2836.3076
$ws.Range("I100").Value = 1621.625
$ws.Range("K100").Value = 1621.625
$ws.Range("M100").Value = -1080.625
$ws.Range("H131").Value = 80000
$ws.Range("I131").Value = 80000
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 80000
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 2203.8293
$ws.Range("I136").Value = 1871.9062
$ws.Range("K136").Value = 5615.7186
$ws.Range("M136").Value = -3065.7186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 6249.75
$ws.Range("I17").Value = 6499.5
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 6499.5
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -6327.5
$ws.Range("N17").Value = -6344
